$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Tracker")

# --- Row 11 (Michael Burton / Documentation row was formerly blank for dates/hours/notes) ---
$ws.Range("D11").Formula = "=TODAY()"
$ws.Range("E11").Formula = "=TODAY()"
$ws.Range("G11").Value = 3
$ws.Range("I11").Value = "Documentation"

# --- Row 12 (Jesse Haynes-Lewis / Documentation row) ---
$ws.Range("D12").Formula = "=TODAY()"
$ws.Range("E12").Formula = "=TODAY()"
$ws.Range("G12").Value = 3
$ws.Range("I12").Value = "Documentation"

# Update the selected cell to reflect where the user finished editing
$ws.Range("P25").Select() | Out-Null
